{"js": "// Apply equation text replacements (two-digit x two-digit multiplication cells).\n// Each \"before\" value is unique within the document, so an exact, case-sensitive\n// search-and-replace on each pair safely targets the correct cell.\nconst replacements = [[\"65\u00d763=\", \"87\u00d794=\"], [\"65\u00d799=\", \"24\u00d779=\"], [\"27\u00d792=\", \"49\u00d789=\"], [\"64\u00d772=\", \"11\u00d735=\"], [\"16\u00d797=\", \"46\u00d736=\"], [\"60\u00d729=\", \"52\u00d715=\"], [\"96\u00d773=\", \"83\u00d739=\"], [\"13\u00d721=\", \"18\u00d711=\"], [\"43\u00d752=\", \"71\u00d758=\"], [\"42\u00d712=\", \"57\u00d739=\"], [\"75\u00d741=\", \"16\u00d762=\"], [\"90\u00d795=\", \"13\u00d754=\"], [\"75\u00d772=\", \"47\u00d733=\"], [\"71\u00d794=\", \"60\u00d757=\"], [\"36\u00d788=\", \"93\u00d753=\"], [\"16\u00d716=\", \"23\u00d751=\"], [\"40\u00d746=\", \"30\u00d722=\"], [\"92\u00d779=\", \"59\u00d772=\"], [\"26\u00d761=\", \"60\u00d774=\"], [\"86\u00d777=\", \"46\u00d712=\"], [\"83\u00d712=\", \"83\u00d713=\"], [\"22\u00d749=\", \"52\u00d793=\"], [\"66\u00d715=\", \"75\u00d788=\"], [\"23\u00d750=\", \"64\u00d785=\"], [\"88\u00d771=\", \"54\u00d720=\"]];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update two-digit x two-digit multiplication equation cells to new values.\n# Each 'before' equation text is unique in the document, so Find/Execute with\n# wdReplaceOne safely targets exactly the right cell for every pair.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"65\u00d763=\", \"87\u00d794=\"),\n    @(\"65\u00d799=\", \"24\u00d779=\"),\n    @(\"27\u00d792=\", \"49\u00d789=\"),\n    @(\"64\u00d772=\", \"11\u00d735=\"),\n    @(\"16\u00d797=\", \"46\u00d736=\"),\n    @(\"60\u00d729=\", \"52\u00d715=\"),\n    @(\"96\u00d773=\", \"83\u00d739=\"),\n    @(\"13\u00d721=\", \"18\u00d711=\"),\n    @(\"43\u00d752=\", \"71\u00d758=\"),\n    @(\"42\u00d712=\", \"57\u00d739=\"),\n    @(\"75\u00d741=\", \"16\u00d762=\"),\n    @(\"90\u00d795=\", \"13\u00d754=\"),\n    @(\"75\u00d772=\", \"47\u00d733=\"),\n    @(\"71\u00d794=\", \"60\u00d757=\"),\n    @(\"36\u00d788=\", \"93\u00d753=\"),\n    @(\"16\u00d716=\", \"23\u00d751=\"),\n    @(\"40\u00d746=\", \"30\u00d722=\"),\n    @(\"92\u00d779=\", \"59\u00d772=\"),\n    @(\"26\u00d761=\", \"60\u00d774=\"),\n    @(\"86\u00d777=\", \"46\u00d712=\"),\n    @(\"83\u00d712=\", \"83\u00d713=\"),\n    @(\"22\u00d749=\", \"52\u00d793=\"),\n    @(\"66\u00d715=\", \"75\u00d788=\"),\n    @(\"23\u00d750=\", \"64\u00d785=\"),\n    @(\"88\u00d771=\", \"54\u00d720=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute(\n        [ref]$oldText,\n        [ref]$true,\n        [ref]$false,\n        [ref]$false,\n        [ref]$false,\n        [ref]$false,\n        [ref]$true,\n        [ref]1,\n        [ref]$false,\n        [ref]$newText,\n        [ref]1   # wdReplaceOne\n    ) | Out-Null\n}\n"}
